$d = $word.ActiveDocument

# 1) Update placeholders in the first "account info" paragraph: add "_1" suffix
#    to Account_Number, Address and Debt_Amount (Status_date stays unchanged).
$d.Content.Find.Execute("{Account_Number}", $true, $false, $false, $false, $false, $true, 1, $false, "{Account_Number_1}", 2) | Out-Null
$d.Content.Find.Execute("{Address}", $true, $false, $false, $false, $false, $true, 1, $false, "{Address_1}", 2) | Out-Null
$d.Content.Find.Execute("{Debt_Amount}", $true, $false, $false, $false, $false, $true, 1, $false, "{Debt_Amount_1}", 2) | Out-Null

# 2) Locate the paragraph that holds the "Будь ласка, не ігноруйте..." warning
#    which currently has no direct paragraph formatting (no <w:pPr>), and give
#    it the same paragraph mark run properties used throughout the letter.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ігноруйте*") {
        $targetPara = $p
    }
}

$targetParaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:wordDocument xmlns:w="http://schemas.microsoft.com/office/word/2003/wordml">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>Будь ласка, не ігноруйте це повідомлення. Несплата заборгованості може призвести до подальших правових наслідків, включаючи передачу справи до суду.</w:t>
      </w:r>
    </w:p>
</w:body>
</w:wordDocument>
"@
$targetPara.Range.InsertXML($targetParaXml)

# 3) Replace the trailing empty paragraph with a second copy of the whole
#    letter (separator line, greeting, account info with "_2" placeholders,
#    reminder text and the closing warning), matching the personalized
#    letter above so the single document now addresses two recipients.
$lastPara = $d.Paragraphs.Last
$newBlockXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:wordDocument xmlns:w="http://schemas.microsoft.com/office/word/2003/wordml">
<w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>----------------------------------------------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>Шановний клієнте!</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>Ми звернули увагу, що за Вашим особовим рахунком № {Account_Number_2} за адресом {Address_2}, накопичилася заборгованість. Станом на {Status_date} сума заборгованості складає {Debt_Amount_2} гривень.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>Ми наполегливо рекомендуємо Вам погасити заборгованість у найкоротші терміни для уникнення нарахування пені та інших санкцій. Ви можете здійснити оплату у будь-якому відділенні банку або через електронні платіжні системи.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:noProof/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="ru-RU" w:eastAsia="ru-RU"/>
        </w:rPr>
        <w:t>Будь ласка, не ігноруйте це повідомлення. Несплата заборгованості може призвести до подальших правових наслідків, включаючи передачу справи до суду.</w:t>
      </w:r>
    </w:p>
</w:body>
</w:wordDocument>
"@
$lastPara.Range.InsertXML($newBlockXml)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
